# Commit: "Fruta / hortaliza, semanal" — add the week's new price record.
#
# The new observation is inserted as row 746 (pushing the former rows
# 746-806 down to 747-807), matching the sheet's reverse-chronological /
# batch-appended layout used throughout this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 746, shifting existing rows 746:806 -> 747:807.
$ws.Rows(746).Insert()

# Populate the new row with the new weekly record.
$ws.Cells.Item(746, 1).Value  = 8
$ws.Cells.Item(746, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(746, 3).Value  = "Coquimbo"
$ws.Cells.Item(746, 4).Value  = (Get-Date -Year 2023 -Month 3 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(746, 5).Value  = 4
$ws.Cells.Item(746, 6).Value  = 100112024
$ws.Cells.Item(746, 7).Value  = "Choclo"
$ws.Cells.Item(746, 8).Value  = "Choclero"
$ws.Cells.Item(746, 9).Value  = "Primera"
$ws.Cells.Item(746, 10).Value = 6000
$ws.Cells.Item(746, 11).Value = 300
$ws.Cells.Item(746, 12).Value = 350
$ws.Cells.Item(746, 13).Value = 325
$ws.Cells.Item(746, 14).Value = "$/unidad"
$ws.Cells.Item(746, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(746, 16).Value = 325
$ws.Cells.Item(746, 17).Value = 1
$ws.Cells.Item(746, 18).Value = "Hortaliza"
